$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add "NA" to column E (duplicate_image_filename) for rows 2 through 21,
# matching the author's commit: "add the NA's under duplicate_image_filename"
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
